$d = $word.ActiveDocument

$replacements = @(
    @("47×56=", "83×59="),
    @("65×89=", "81×76="),
    @("41×13=", "91×60="),
    @("20×48=", "31×34="),
    @("57×41=", "96×15="),
    @("88×85=", "39×80="),
    @("25×13=", "49×33="),
    @("99×53=", "36×79="),
    @("59×79=", "37×48="),
    @("95×68=", "17×14="),
    @("64×13=", "35×98="),
    @("44×93=", "27×13="),
    @("16×49=", "61×57="),
    @("45×77=", "83×92="),
    @("56×30=", "86×57="),
    @("68×96=", "76×82="),
    @("63×63=", "42×68="),
    @("14×53=", "29×27="),
    @("71×59=", "77×70="),
    @("41×18=", "62×36="),
    @("29×62=", "32×17="),
    @("87×83=", "82×35="),
    @("90×57=", "75×36="),
    @("27×59=", "58×47="),
    @("52×14=", "82×41=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
